# Improve waffle plot coloring and ordering
$wb = $excel.ActiveWorkbook

# --- 1) "Color" sheet: bump saturation 40 -> 50 in the hsva() helper formulas (C2:N6) ---
$wsColor = $wb.Worksheets.Item("Color")
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N")
for ($r = 2; $r -le 6; $r++) {
    foreach ($col in $cols) {
        $cell = $col + $r
        $formula = '="hsva("&' + $col + '$1&",50,"&$B' + $r + '&",100)"'
        $wsColor.Range($cell).Formula = $formula
    }
}

# --- 2) "ColorMap" sheet: re-derive the hue/key rows (reordered + recolored waffle legend) ---
$wsColorMap = $wb.Worksheets.Item("ColorMap")

# Clear out the old A2:B25 block first (row count shrinks from 25 to 15).
$wsColorMap.Range("A2:B25").Clear()

$hues = @(180,150,180,150,180,210,240,120,150,180,210,240,270,300)
$keys = @("main","aux","Total","Capital","State","Population","Significance","Southeast","Northeast","Northwest","Southwest","Steel Belt","Plains","Oconus")
for ($i = 0; $i -lt $hues.Length; $i++) {
    $row = $i + 2
    $wsColorMap.Range("A" + $row).Value = $hues[$i]
    $wsColorMap.Range("B" + $row).Value = $keys[$i]
}

# --- 3) Sheet view / selection bookkeeping (matches the author's final click-state) ---
$wsColor.Activate()
$wsColor.Range("F11").Select()

$wsColorMap.Activate()
$wsColorMap.Range("C10").Select()

$wsCities = $wb.Worksheets.Item("Cities")
$wsCities.Activate()
$wsCities.Range("C24").Select()
